# Minor graphical error correction
# Nudge four "Lookup" round-rect textbox callouts (and the related
# "User Authorization" box) that were mis-positioned on slide 1, moving
# them to their corrected locations.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# TextBox 868 ("Provider Lookup *") - id 869
$shp869 = $s.Shapes.Item(34)
$shp869.Top = 676.8484496968504

# TextBox 938 ("User Bank Lookup *") - id 939
$shp939 = $s.Shapes.Item(66)
$shp939.Left = 174.844131488189
$shp939.Top = 361.2092590984252

# TextBox 942 ("Merchant Lookup *") - id 943
$shp943 = $s.Shapes.Item(70)
$shp943.Top = 613.448547437008

# TextBox 943 ("User Authorization ...") - id 944
$shp944 = $s.Shapes.Item(71)
$shp944.Top = 717.6729736858267
